# Fix an issue with the case the program is called without any parameters.
#
# 1) Slide 10, "Textfeld 22": the run containing a lone space between the
#    bold "bm" run and the "= float [0.0...1.0]" run is removed; the space
#    becomes the leading character of the "= float [0.0...1.0]" run instead
#    (the run's own formatting - dirty="0" - is kept).
# 2) Slide 8: the right-brace shape is resized/repositioned very slightly
#    and the small rectangle labelled with "no of chars" value is moved
#    to sit under the brace's new position.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 10 text fix
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$shpText = $slide10.Shapes.Item(2)          # "Textfeld 22"
$tr = $shpText.TextFrame.TextRange

# Paragraph 19 reads: "Supported parameters: only bm = float [0.0...1.0]"
# (NB: searching/printing text read back through TextRange.Text mangles
# non-ASCII characters such as the ellipsis, so only ASCII needles are used
# to locate offsets; the ellipsis itself is only ever *written*, via the
# literal character below, never compared against read-back text.)
$para = $tr.Paragraphs(19, 1)
$paraText = $para.Text
$bmIdx = $paraText.IndexOf("bm")
$eqIdx = $paraText.IndexOf("= float")
$eqLen = $paraText.Length - $eqIdx

$spaceStart = $para.Start + $bmIdx + 2      # position right after "bm"
$eqStart = $para.Start + $eqIdx             # position where "=" begins

# Preserve the shape's current (auto-fitted) height, because re-assigning
# any text in this autosized textbox makes the host recompute the fit.
$origHeight = $shpText.Height

# Prepend the space onto the "= float ..." run, keeping that run's own
# formatting intact.
$eqRange = $tr.Characters($eqStart, $eqLen)
$eqRange.Text = " = float [0.0" + [char]0x2026 + "1.0]"

# Remove the now-redundant standalone space run.
$spaceRange = $tr.Characters($spaceStart, 1)
$spaceRange.Text = ""

# Restore the shape's height (spAutoFit is left untouched).
$shpText.Height = $origHeight

# ---------------------------------------------------------------------------
# 2) Slide 8 shape geometry fix
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

$brace = $slide8.Shapes.Item(24)            # "Geschweifte Klammer rechts 36"
$brace.Top = 98.4747314453125
$brace.Height = 360.56341552734375

$rect = $slide8.Shapes.Item(29)             # "Rechteck 51"
$rect.Left = 470.7465515136719
